$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 4)
$ws.Range("V4").Value = "Milestone"
$ws.Range("W4").Value = "TTO"
$ws.Range("X4").Value = "ADMC"

# New data cells (row 5)
$ws.Range("V5").Value = "Player Controller"
$ws.Range("X5").Value = -4

# Column width for V (used to display "Player Controller" nicely)
# (closest achievable value to the target 18.5546875 given the engine's
# internal 1/6-character rounding of ColumnWidth)
$ws.Range("V1").ColumnWidth = 17.6

# Update view: scroll so column B is the left-most visible column,
# and set the active selection to V7
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("V7").Select()
